$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab per the workbook.xml sheet name change.
$ws.Name = "o554F"

# Tiny floating point recalculation adjustments in existing rows (row 13 and row 15).
$ws.Range("E13").Value = 0.9828962381417687
$ws.Range("G13").Value = 0.9876195213177464
$ws.Range("K13").Value = 0.9859675287895707
$ws.Range("L13").Value = 0.994105520707096

$ws.Range("D15").Value = 1.403486764409819
$ws.Range("I15").Value = 0.9211748398454243
$ws.Range("J15").Value = 1.403486764409819

# New row 16 of data.
$ws.Range("A16").Value = 14
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats - reuse row 15's bold/border style
$excel.CutCopyMode = $false

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.749800258207028
$ws.Range("D16").Value = 0.2761562175175173
$ws.Range("E16").Value = 1.047988804147873
$ws.Range("F16").Value = 1.749800258207028
$ws.Range("G16").Value = 0.5946063655581701
$ws.Range("H16").Value = 1.16863442246823
$ws.Range("I16").Value = 1.13665739027051
$ws.Range("J16").Value = 0.2761562175175173
$ws.Range("K16").Value = 0.6620725108326952
$ws.Range("L16").Value = 1.205936384519862
$ws.Range("M16").Value = 0.9956405763615549
